$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold text on dark fill, same as existing CHANNEL CODE/
# CHANNEL NAME headers in E2:F2) onto the new header cells before writing
# their text, so the shared style index is preserved on save.
$ws.Range("E2").Copy()
$ws.Range("G2:J2").PasteSpecial(-4122)

# New header row values (PROVINCE, CITY/TOWN, BARANGAY, STREET)
$ws.Range("G2").Value = "PROVINCE"
$ws.Range("H2").Value = "CITY/TOWN"
$ws.Range("I2").Value = "BARANGAY"
$ws.Range("J2").Value = "STREET"

# New sample data row values (MANILA, MAKATI, SAN ANTONIO, ARANGA)
$ws.Range("G3").Value = "MANILA"
$ws.Range("H3").Value = "MAKATI"
$ws.Range("I3").Value = "SAN ANTONIO"
$ws.Range("J3").Value = "ARANGA"

# Move the active selection to B3, matching the saved workbook view state.
$ws.Range("B3").Select()
